$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.980.37'
$ws.Range('E2').Value = '  +0.42%  '

$ws.Range('D3').Value = '1.884.18'
$ws.Range('E3').Value = '  -0.17%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').Value = '  -0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7450'
$ws.Range('E5').Value = '  -2.81%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.90'
$ws.Range('E6').Value = '  +0.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9995'
$ws.Range('E7').Value = '  -0.16%  '

$ws.Range('E8').Value = '  +1.13%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07237'
$ws.Range('E9').Value = '  +1.62%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.91'
$ws.Range('E10').Value = '  -2.64%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08348'
$ws.Range('E11').Value = '  -2.16%  '

$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7564'
$ws.Range('E12').Value = '  -0.83%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.423'
$ws.Range('E13').Value = '  +1.15%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.816.06'
$ws.Range('E14').Value = '  -5.29%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.57'
$ws.Range('E15').Value = '  -1.01%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.140'
$ws.Range('E16').Value = '  +0.14%  '

$ws.Range('D17').Value = '29.910.74'
$ws.Range('E17').Value = '  -0.17%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '249.97'
$ws.Range('E18').Value = '  +2.40%  '

$ws.Range('E19').Value = '  -0.96%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007856'
$ws.Range('E20').Value = '  +0.50%  '

$ws.Range('D21').Value = '2.195.12'
$ws.Range('E21').Value = '  +0.99%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9980'
$ws.Range('E22').Value = '  -0.17%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.998'
$ws.Range('E23').Value = '  -0.19%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  -0.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1567'
$ws.Range('E25').Value = '  -3.38%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.291'
$ws.Range('E26').Value = '  -1.05%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.47'
$ws.Range('E27').Value = '  +1.38%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.72'
$ws.Range('E28').Value = '  -0.27%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.042'
$ws.Range('E29').Value = '  +0.42%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.487'
$ws.Range('E30').Value = '  -1.49%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.620'
$ws.Range('E31').Value = '  +2.48%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.536'
$ws.Range('E32').Value = '  -0.28%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.233'
$ws.Range('E33').Value = '  +2.69%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05380'
$ws.Range('E34').Value = '  -1.06%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.256'
$ws.Range('E35').Value = '  +1.09%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7558'
$ws.Range('E36').Value = '  +1.43%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9934'
$ws.Range('E37').Value = '  -0.78%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.707'
$ws.Range('E38').Value = '  +0.21%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01964'
$ws.Range('E39').Value = '  +0.93%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.762'
$ws.Range('E40').Value = '  -0.71%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4558'
$ws.Range('E41').Value = '  +2.12%  '

$ws.Range('D42').Value = '1.105.48'
$ws.Range('E42').Value = '  +0.42%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.039'
$ws.Range('E43').Value = '  -0.69%  '

$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.58'
$ws.Range('E44').Value = '  -0.57%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8688'
$ws.Range('E45').Value = '  +1.57%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '104.51'
$ws.Range('E46').Value = '  +1.45%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.000'
$ws.Range('E47').Value = '  -0.04%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.865'
$ws.Range('E48').Value = '  -0.22%  '

$ws.Range('E49').Value = '  -0.65%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.569'
$ws.Range('E50').Value = '  -1.65%  '

$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.039.07'
$ws.Range('E51').Value = '  -0.08%  '
